$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("Z1").Style = "Normal"
$ws.Range("Z1").Value = "x"
$ws.Range("Z1").ClearContents()
